$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column AG (16-jul) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy formatting from the previous day's column (AF) into the new column (AG)
$wsPrix.Range("AF1:AF25").Copy()
$wsPrix.Range("AG1:AG25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsPrix.Range("AG1").Value = "16-jul"
$wsPrix.Range("AG2").Value = 74.97
$wsPrix.Range("AG3").Value = 56.53
$wsPrix.Range("AG4").Value = 55.26
$wsPrix.Range("AG5").Value = 48.19
$wsPrix.Range("AG6").Value = 44.71
$wsPrix.Range("AG7").Value = 46.97
$wsPrix.Range("AG8").Value = 52.33
$wsPrix.Range("AG9").Value = 62.55
$wsPrix.Range("AG10").Value = 75.2
$wsPrix.Range("AG11").Value = 66.98999999999999
$wsPrix.Range("AG12").Value = 42.63
$wsPrix.Range("AG13").Value = 28.81
$wsPrix.Range("AG14").Value = 53.6
$wsPrix.Range("AG15").Value = 36.26
$wsPrix.Range("AG16").Value = 33.25
$wsPrix.Range("AG17").Value = 29.01
$wsPrix.Range("AG18").Value = 39.76
$wsPrix.Range("AG19").Value = 57.57
$wsPrix.Range("AG20").Value = 79.09999999999999
$wsPrix.Range("AG21").Value = 109.62
$wsPrix.Range("AG22").Value = 122.84
$wsPrix.Range("AG23").Value = 108.4
$wsPrix.Range("AG24").Value = 117.4
$wsPrix.Range("AG25").Value = 111.14

# --- Sheet "Gaz": append row 30 (2025-07-14) ---
# The date-like string is entered via a text formula and then converted to a
# static value in place, so Excel does not auto-convert it into a date serial
# number (which would happen with a plain .Value assignment).
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A30").Formula = "=""2025-07-14"""
$wsGaz.Range("A30").Copy()
$wsGaz.Range("A30").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$wsGaz.Range("B30").Value = 34.275

# --- Sheet "CO2": append row 30 (2025-07-14) ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A30").Formula = "=""2025-07-14"""
$wsCO2.Range("A30").Copy()
$wsCO2.Range("A30").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$wsCO2.Range("B30").Value = 69.59999999999999

Write-Output "edit complete"
